# Phase 1 development done
# Add a new bullet-style tip after the "Be careful with armor!" paragraph:
#   " - Use your spaceship! Can be used to store items, in case you die"

$d = $word.ActiveDocument

# The new tip goes at the very end of the document, right after the last
# paragraph ("Be careful with armor! ..."). Grab that last paragraph and
# collapse its range to its end point so we can append a new paragraph mark.
$lastPara = $d.Paragraphs.Last
$endRange = $lastPara.Range
$endRange.Collapse(0)              # wdCollapseEnd
$endRange.InsertParagraphAfter()   # creates a brand-new (default-formatted) paragraph

# The newly created paragraph is now the last paragraph in the document.
# Insert the tip text into it.
$newPara = $d.Paragraphs.Last
$newRange = $newPara.Range
$newRange.Collapse(0)              # wdCollapseEnd
$newRange.InsertAfter(" - Use your spaceship! Can be used to store items, in case you die")
